# Juno: check in to OLPRODLOC.
#
# Bold the lead-in word(s) of the intro paragraph and each top-level
# "Dies umfasst:" list item, turning them into "Label: sentence" style
# bullets, and refresh the German wording of each sentence.

$d = $word.ActiveDocument

function Apply-Edit {
    param(
        [int]$ParaIndex,
        [string]$OldText,
        [string]$BoldPrefix,
        [string]$Rest,
        [int]$RestColor
    )

    $para = $d.Paragraphs($ParaIndex).Range

    $newText = $BoldPrefix + $Rest
    $found = $para.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Host "NOT FOUND (para" $ParaIndex "):" $OldText
        return
    }

    $paraStart = $d.Paragraphs($ParaIndex).Range.Start

    # Bold just the lead-in label.
    $boldRange = $d.Range($paraStart, $paraStart + $BoldPrefix.Length)
    $boldRange.Bold = 1

    # Make sure the remainder is explicitly not bold, and recolor it if needed.
    $restRange = $d.Range($paraStart + $BoldPrefix.Length, $paraStart + $newText.Length)
    $restRange.Bold = 0
    if ($RestColor -ge 0) {
        $restRange.Font.Color = $RestColor
    }
}

Apply-Edit 3 `
    "Dieses Dokument umreißt die Hauptverantwortlichkeiten aller Mitglieder des Designteams am Graphic Design Institute." `
    "Zweck" `
    ": Dieses Dokument beschreibt die Kernaufgaben aller Designteammitglieder am Graphic Design Institute." `
    -1

Apply-Edit 5 `
    "Zusammenarbeit mit anderen Designerinnen und Designern, Entwickelnden und Beteiligten, um qualitativ hochwertige Designs zu erstellen, die den Projektanforderungen entsprechen." `
    "Zusammenarbeit: Arbeiten Sie" `
    " gemeinsam mit anderen Designern, Entwicklern und Projektbeteiligten zusammen, um hochwertige Designs zu erstellen, die den Projektanforderungen entsprechen." `
    -1

Apply-Edit 13 `
    "Erstellung visuell ansprechender Designs, die benutzerfreundlich, zugänglich und bedarfsgerecht sind." `
    "Design" `
    ": Erstellen Sie visuell ansprechende Designs, die benutzerfreundliche, barrierefrei und reaktionsfähig sind." `
    -1

Apply-Edit 22 `
    "Effektive Kommunikation mit Teammitgliedern, Beteiligten und der Kundschaft, um sicherzustellen, dass die Projektanforderungen erfüllt werden." `
    "Kommunikation" `
    ": Kommunizieren Sie effektiv mit Teammitgliedern, Projektbeteiligten und Kunden, um sicherzustellen, dass die Projektanforderungen erfüllt sind." `
    -1

Apply-Edit 30 `
    "Durchführung von Recherchen zur Ermittlung der Bedürfnisse, Präferenzen und Verhaltensweisen von Benutzerinnen und Benutzern als Grundlage für Designentscheidungen." `
    "Forschung" `
    ": Führen Sie Untersuchungen durch, um Benutzerbedürfnisse, Vorlieben und Verhaltensweisen zu identifizieren, um Designentscheidungen zu treffen." `
    -1

Apply-Edit 39 `
    "Durchführung von Gebrauchstauglichkeitsprüfungen, um sicherzustellen, dass die Entwürfe den Bedürfnissen der Zielgruppe entsprechen und für alle Benutzerinnen und Benutzern zugänglich sind." `
    "Tests" `
    ": Führen Sie Benutzerfreundlichkeitstests durch, um sicherzustellen, dass Designs den Anforderungen der Benutzer entsprechen und für alle Benutzer zugänglich sind." `
    -1

Apply-Edit 48 `
    "Erstellung und Pflege von Designdokumentation, einschließlich Designspezifikationen, Stilanleitungen und Entwurfsmuster." `
    "Dokumentation" `
    ": Erstellen und Verwalten von Entwurfsdokumentationen, einschließlich Entwurfsspezifikationen, Stilführungslinien und Entwurfsmustern." `
    -1

Apply-Edit 56 `
    "Sich über die neuesten Designtrend, Tools und Technologien auf dem Laufenden halten, um die Designqualität und -effizienz zu verbessern." `
    "Professionelle Entwicklung" `
    ": Bleiben Sie mit den neuesten Designtrends, Tools und Technologien auf dem laufenden, um die Designqualität und -effizienz zu verbessern." `
    -1

Apply-Edit 61 `
    "Leitung des Designteams und Anleitung von Design-Nachwuchskräften." `
    "Führung:" `
    " Führen Sie das Designteam und bieten Ihnen Anleitungen für Juniordesigner." `
    1118481

Write-Host "Done."
